# Update R script and sample file
# Adds a new "Sheet3" at the end of the workbook containing two FST matrices
# (mean and weighted Weir & Cockerham FST) pasted in from an R analysis.

$wb = $excel.ActiveWorkbook

# --- Add Sheet3 after the last existing sheet --------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet3"

$pops = @("SU18", "AR18", "SI18", "SD18", "YS21", "IS21", "UH21")

$table1 = @(
    @("SU18", "NA", "NA", "NA", "NA", "NA", "NA"),
    @("AR18", 0.134, "NA", "NA", "NA", "NA", "NA"),
    @("SI18", 0.007, 0.049, "NA", "NA", "NA", "NA"),
    @("SD18", 0.106, 0.102, 0.027, "NA", "NA", "NA"),
    @("YS21", 0.011, 0.029, 0.468, 0.102, "NA", "NA"),
    @("IS21", 0.004, 0.015, 0.042, 0.083, 0.374, "NA"),
    @("UH21", 0, 0, 0, 0, 0, 0)
)

$table2 = @(
    @("SU18", "NA", "NA", "NA", "NA", "NA", "NA"),
    @("AR18", 0.002749224, "NA", "NA", "NA", "NA", "NA"),
    @("SI18", 0.00783212, 0.004766385, "NA", "NA", "NA", "NA"),
    @("SD18", 0.003300563, 0.004062428, 0.0059278723, "NA", "NA", "NA"),
    @("YS21", 0.005832439, 0.004269984, 0.0002982616, 0.003294308, "NA", "NA"),
    @("IS21", 0.006783323, 0.005118175, 0.0045439467, 0.004380776, 0.0009208749, "NA"),
    @("UH21", 0.045215666, 0.032584532, 0.0213079975, 0.041324086, 0.0133642404, 0.03022955)
)

# The original author pasted these matrices in such a way that the new
# shared-string pool picks up "NA" before "POP"/"Pop" -- fill the numeric /
# NA cell bodies first, then the row labels, then the two header rows, so the
# workbook's shared-string table is built in the same order.

for ($r = 0; $r -lt $table1.Length; $r++) {
    $row = $table1[$r]
    $excelRow = $r + 2
    for ($c = 1; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
    $ws.Cells.Item($excelRow, 8).Value = "NA"
}

for ($r = 0; $r -lt $table2.Length; $r++) {
    $row = $table2[$r]
    $excelRow = $r + 12
    for ($c = 1; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
    $ws.Cells.Item($excelRow, 8).Value = "NA"
}

# Row labels (column A), reusing population-name strings already in the pool
for ($r = 0; $r -lt $table1.Length; $r++) {
    $ws.Cells.Item($r + 2, 1).Value = $table1[$r][0]
}
for ($r = 0; $r -lt $table2.Length; $r++) {
    $ws.Cells.Item($r + 12, 1).Value = $table2[$r][0]
}

# Header rows: A1 = "POP", A11 = "Pop"; B:H reuse the population-name strings
$ws.Cells.Item(1, 1).Value = "POP"
for ($c = 0; $c -lt $pops.Length; $c++) {
    $ws.Cells.Item(1, $c + 2).Value = $pops[$c]
}

$ws.Cells.Item(11, 1).Value = "Pop"
for ($c = 0; $c -lt $pops.Length; $c++) {
    $ws.Cells.Item(11, $c + 2).Value = $pops[$c]
}

# --- Page setup (matches Print Layout the author left set on the new sheet) --
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection: the author had the second table selected when saving --------
$ws.Range("A11:H18").Select()
